# The "hasInterval" properties row (row 17: hasInterval / Time interval /
# Zeitintervall / hasSequenceBounds / IntervalValue / Interval) was a
# left-over / incorrect entry. Select the whole row the way a user would
# before removing it, then delete it so everything below shifts up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$null = $ws.Rows.Item(17).Select()
$ws.Rows.Item(17).Delete()
